# Fixed epoch removal code: the previous version of this sheet accidentally
# dropped two genuine AGN outlier rows when the stray "111353.73+515725.8 "
# epoch row got cleaned up, and it was also missing a repeat-epoch reading
# for one of the W1 sources. Put all of that back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row after the existing "121449.54+572734.1 / W1 / 3" row
#    (currently row 14) with a second Epoch reading for that AGN/band. Do
#    this first so it doesn't shift the two brand-new rows appended below.
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "121449.54+572734.1"
$ws.Range("B15").Value = "W1"
$ws.Range("C15").Value = 4

# 2) Append the first of the two previously-removed AGN rows.
$ws.Range("A29").Value = "141648.88+530903.5"
$ws.Range("B29").Value = "W1"
$ws.Range("C29").Value = 23

# 3) Trim the stray trailing space on the "111353.73+515725.8 " shared string
#    (rows 5 and 6 both use it).
$ws.Range("A5").Value = "111353.73+515725.8"
$ws.Range("A6").Value = "111353.73+515725.8"

# 4) Append the second previously-removed AGN row.
$ws.Range("A30").Value = "125731.87+272313.3"
$ws.Range("B30").Value = "W1"
$ws.Range("C30").Value = 6

# 5) Restore the on-screen selection to match the saved view state.
$ws.Range("G24").Select()
